$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-apply the built-in "No Style, No Grid" table style to the three
#    tables in the deck (slides 14, 15, 16), replacing the previous
#    custom "Table_0" style.
# ---------------------------------------------------------------------
$newTableStyleId = "{67F619A2-B602-402B-A6B4-33E8EC55AE31}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Switch the presentation's design theme colors from the custom
#    "Integral" (Red Violet) palette to the default "Office Theme"
#    palette.
# ---------------------------------------------------------------------
function Get-BgrLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = Get-BgrLong $officeThemeColors[$i - 1]
}
